$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header label in H1
$ws.Range("H1").Value = "measured"

# For rows 2-5, copy the old A:F values (the "design" values) into H:M,
# preserving the original (pre-edit) wall thickness in column M,
# then shrink the wall thickness in column F by an eighth (store the
# newly measured/reduced values).

$rows = 2, 3, 4, 5
foreach ($r in $rows) {
    $ws.Range("H$r").Value = $ws.Range("A$r").Value2
    $ws.Range("I$r").Value = $ws.Range("B$r").Value2
    $ws.Range("J$r").Value = $ws.Range("C$r").Value2
    $ws.Range("K$r").Value = $ws.Range("D$r").Value2
    $ws.Range("L$r").Value = $ws.Range("E$r").Value2
    $ws.Range("M$r").Value = $ws.Range("F$r").Value2
}

# Decrease the wall thickness (column F) by an eighth for these rows
$ws.Range("F2").Value = 0.35
$ws.Range("F3").Value = 0.25
$ws.Range("F4").Value = 0.12
$ws.Range("F5").Value = 0.07

# Update the selected cell to match the author's final cursor position
$ws.Range("I7").Select()

$wb.Save()
